# Apply the edit described by the diff:
# - Remove the first 10 data rows (old rows 2-12), shifting remaining
#   data rows (old 13-22) up to become new rows 2-11.
# - Append 10 new data rows (new rows 12-21) with freshly generated values.
# - Net effect: sheet dimension goes from A1:C22 to A1:C21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old rows 2 through 12 (11 rows), which shifts rows 13-22 up.
$ws.Range("A2:C12").EntireRow.Delete() | Out-Null

# New data to append at the bottom (rows 12-21 after the shift above).
$newRows = @(
    @(-0.3932898044586127, -2.516493201255787, -11.63133525848387),
    @(-0.3210607767105248, 2.201687335968037, -5.495597958564743),
    @(-3.88015073537826, 7.319447636604309, -2.401085853576681),
    @(-2.474413871765147, 7.422795295715336, -7.479803562164315),
    @(-4.692895889282228, 8.487199664115906, -9.337096989154816),
    @(-3.639542102813719, 8.090452075004578, -3.927529096603386),
    @(-5.387722790241247, 3.111244738101949, 8.269636750221276),
    @(-1.612288236618022, 0.4680981636047372, 12.90354442596435),
    @(0.8442984223365615, -0.4397069215774612, 9.103015005588517),
    @(-3.634706258773804, -4.129897594451904, 6.328503251075745)
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
